$d = $word.ActiveDocument

$d.Content.Find.Execute("155÷3=51, 2", $true, $false, $false, $false, $false, $true, 1, $false, "762÷7=108, 6", 2) | Out-Null
$d.Content.Find.Execute("345÷5=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "543÷4=135, 3", 2) | Out-Null
$d.Content.Find.Execute("473÷3=157, 2", $true, $false, $false, $false, $false, $true, 1, $false, "535÷7=76, 3", 2) | Out-Null
$d.Content.Find.Execute("365÷6=60, 5", $true, $false, $false, $false, $false, $true, 1, $false, "368÷9=40, 8", 2) | Out-Null
$d.Content.Find.Execute("972÷8=121, 4", $true, $false, $false, $false, $false, $true, 1, $false, "253÷5=50, 3", 2) | Out-Null
$d.Content.Find.Execute("253÷3=84, 1", $true, $false, $false, $false, $false, $true, 1, $false, "324÷2=162, 0", 2) | Out-Null
$d.Content.Find.Execute("489÷4=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "252÷4=63, 0", 2) | Out-Null
$d.Content.Find.Execute("653÷8=81, 5", $true, $false, $false, $false, $false, $true, 1, $false, "229÷3=76, 1", 2) | Out-Null
$d.Content.Find.Execute("896÷7=128, 0", $true, $false, $false, $false, $false, $true, 1, $false, "222÷8=27, 6", 2) | Out-Null
$d.Content.Find.Execute("186÷5=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "417÷9=46, 3", 2) | Out-Null
$d.Content.Find.Execute("692÷4=173, 0", $true, $false, $false, $false, $false, $true, 1, $false, "565÷2=282, 1", 2) | Out-Null
$d.Content.Find.Execute("416÷5=83, 1", $true, $false, $false, $false, $false, $true, 1, $false, "176÷2=88, 0", 2) | Out-Null
$d.Content.Find.Execute("420÷8=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "415÷8=51, 7", 2) | Out-Null
$d.Content.Find.Execute("278÷9=30, 8", $true, $false, $false, $false, $false, $true, 1, $false, "295÷6=49, 1", 2) | Out-Null
$d.Content.Find.Execute("884÷5=176, 4", $true, $false, $false, $false, $false, $true, 1, $false, "908÷5=181, 3", 2) | Out-Null
$d.Content.Find.Execute("978÷6=163, 0", $true, $false, $false, $false, $false, $true, 1, $false, "742÷3=247, 1", 2) | Out-Null
$d.Content.Find.Execute("178÷9=19, 7", $true, $false, $false, $false, $false, $true, 1, $false, "986÷7=140, 6", 2) | Out-Null
$d.Content.Find.Execute("750÷7=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "435÷4=108, 3", 2) | Out-Null
$d.Content.Find.Execute("351÷2=175, 1", $true, $false, $false, $false, $false, $true, 1, $false, "881÷7=125, 6", 2) | Out-Null
$d.Content.Find.Execute("868÷6=144, 4", $true, $false, $false, $false, $false, $true, 1, $false, "966÷6=161, 0", 2) | Out-Null
$d.Content.Find.Execute("800÷6=133, 2", $true, $false, $false, $false, $false, $true, 1, $false, "631÷3=210, 1", 2) | Out-Null
$d.Content.Find.Execute("354÷4=88, 2", $true, $false, $false, $false, $false, $true, 1, $false, "694÷3=231, 1", 2) | Out-Null
$d.Content.Find.Execute("581÷9=64, 5", $true, $false, $false, $false, $false, $true, 1, $false, "606÷3=202, 0", 2) | Out-Null
$d.Content.Find.Execute("269÷9=29, 8", $true, $false, $false, $false, $false, $true, 1, $false, "504÷9=56, 0", 2) | Out-Null
$d.Content.Find.Execute("462÷5=92, 2", $true, $false, $false, $false, $false, $true, 1, $false, "450÷7=64, 2", 2) | Out-Null
